$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48/49 swap: Aave and BabyDogeCoin change order, with new D/E values
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'65.32"
$ws.Range("E48").Value = "  +2.45%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000122"
$ws.Range("E49").Value = "  -1.55%  "

# Price (D) and Volume(1h) (E) updates for the rest of the rows
$ws.Range("D2").Value = "'29.228.89"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "'1.844.46"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'242.79"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'44.95"
$ws.Range("E8").Value = "  +6.54%  "
$ws.Range("D9").Value = "'0.07458"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.2963"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "'23.40"
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").Value = "'0.07751"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'1.856.68"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "'5.028"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "'0.6752"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "'83.29"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").Value = "'6.182"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'0.000008772"
$ws.Range("E18").Value = "  +5.32%  "
$ws.Range("D19").Value = "'29.142.33"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "'2.093.05"
$ws.Range("E20").Value = "  +3.93%  "
$ws.Range("D21").Value = "'12.57"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'227.48"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "'0.9997"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'7.186"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'0.9994"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'158.13"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Value = "'8.642"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "'0.1403"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'18.07"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'1.510"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "'4.140"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").Value = "'4.056"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "'1.195"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").Value = "'0.05368"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").Value = "'1.858"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("D36").Value = "'0.7469"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "'1.160"
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("D38").Value = "'2.644"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").Value = "'1.303.99"
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("D40").Value = "'0.01798"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").Value = "'2.759"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").Value = "'6.418"
$ws.Range("E42").Value = "  +7.83%  "
$ws.Range("D43").Value = "'0.9090"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "'0.9992"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'103.50"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'0.08073"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").Value = "'1.991.36"
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("D50").Value = "'0.5135"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").Value = "'1.753"
$ws.Range("E51").Value = "  -0.89%  "
